# Fixed naive component forecaster bug - Presentation state 11.02.
#
# Updates the y_0_forecast (C) and y_1_forecast (E) columns with corrected
# values, and removes the erroneous C2, E2 and C3 cells entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray forecast values that should not exist on rows 2 and 3.
$ws.Cells.Item(2, 3).ClearContents()   # C2
$ws.Cells.Item(2, 5).ClearContents()   # E2
$ws.Cells.Item(3, 3).ClearContents()   # C3

# Corrected values for the remaining rows.
$ws.Cells.Item(3, 5).Value  = 2.957479223435744    # E3

$ws.Cells.Item(4, 3).Value  = -0.01587181126745385 # C4
$ws.Cells.Item(4, 5).Value  = 3.8351443707757      # E4

$ws.Cells.Item(5, 3).Value  = -0.02256889165886955 # C5
$ws.Cells.Item(5, 5).Value  = -0.6757980944263275  # E5

$ws.Cells.Item(6, 3).Value  = 0.09611428386595566  # C6

$ws.Cells.Item(7, 5).Value  = -0.971238541762387   # E7

$ws.Cells.Item(8, 3).Value  = -0.001350220946472191 # C8
$ws.Cells.Item(8, 5).Value  = 0.6008487920565075    # E8

$ws.Cells.Item(9, 5).Value  = -1.58998093318411    # E9

$ws.Cells.Item(10, 3).Value = -0.5761528471665334  # C10
$ws.Cells.Item(10, 5).Value = 0.4501721032283301   # E10

$ws.Cells.Item(11, 5).Value = -0.150175137493469   # E11

$ws.Cells.Item(13, 5).Value = 2.372078088364704    # E13

$ws.Cells.Item(14, 3).Value = -0.4278219446121501  # C14
$ws.Cells.Item(14, 5).Value = -2.378564786744752   # E14

$ws.Cells.Item(15, 3).Value = -1.026566979837429   # C15

$ws.Cells.Item(17, 3).Value = 0.4636049209196802   # C17

$ws.Cells.Item(18, 3).Value = 0.6216390921348403   # C18
$ws.Cells.Item(18, 5).Value = -1.097580983230539   # E18

$ws.Cells.Item(19, 3).Value = -0.6768900623516871  # C19

$wb.Save()
